$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-FormatOnly {
    param($srcRange, $dstRange)
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# Template cells that already carry the two "class group" highlight looks we
# need to replicate onto newly-categorised rows:
#   family "Solo"  -> single-member class highlight (A/B/C, D, E)
#   family "Group" -> multi-member class highlight   (A/B/C, D, E)
$tmplSoloABC = $ws.Range("A49")
$tmplSoloD   = $ws.Range("D49")
$tmplSoloE   = $ws.Range("E49")

$tmplGroupABC = $ws.Range("A4")
$tmplGroupD   = $ws.Range("D4")
$tmplGroupE   = $ws.Range("E4")

function Fill-RowSolo {
    param([int]$row, [string]$classVal, [string]$comment)
    Copy-FormatOnly $tmplSoloABC $ws.Range("A$row")
    Copy-FormatOnly $tmplSoloABC $ws.Range("B$row")
    Copy-FormatOnly $tmplSoloABC $ws.Range("C$row")
    Copy-FormatOnly $tmplSoloD   $ws.Range("D$row")
    Copy-FormatOnly $tmplSoloE   $ws.Range("E$row")
    $ws.Range("A$row").Value2 = $classVal
    $ws.Range("B$row").Value2 = "Done"
    if ($comment) { $ws.Range("C$row").Value2 = $comment }
}

function Fill-RowGroup {
    param([int]$row, [string]$classVal, [string]$comment)
    Copy-FormatOnly $tmplGroupABC $ws.Range("A$row")
    Copy-FormatOnly $tmplGroupABC $ws.Range("B$row")
    Copy-FormatOnly $tmplGroupABC $ws.Range("C$row")
    Copy-FormatOnly $tmplGroupD   $ws.Range("D$row")
    Copy-FormatOnly $tmplGroupE   $ws.Range("E$row")
    $ws.Range("A$row").Value2 = $classVal
    $ws.Range("B$row").Value2 = "Done"
    if ($comment) { $ws.Range("C$row").Value2 = $comment }
}

# --- row 72: newly finished "TTableLayer" function joins the sheet ---
Fill-RowSolo 72 "TTableLayer" $null

# --- row 163: newly finished "control" function (still a lone member here) ---
Fill-RowSolo 163 "control" $null

# --- row 276: another "control" function -> now a multi-member group ---
Fill-RowGroup 276 "control" $null

# --- row 412 (TCircle) already categorised; just mark it Done ---
$ws.Range("B412").Value2 = "Done"

# --- brand-new "TEdgeBox" class: 3 functions done ---
Fill-RowGroup 439 "TEdgeBox" "not needed"
Fill-RowGroup 440 "TEdgeBox" $null
Fill-RowGroup 441 "TEdgeBox" $null

# --- brand-new "TEdgeManager" class: 12 functions done ---
Fill-RowSolo 442 "TEdgeManager" "not needed"
Fill-RowSolo 443 "TEdgeManager" $null
Fill-RowSolo 444 "TEdgeManager" $null
Fill-RowSolo 445 "TEdgeManager" $null
Fill-RowSolo 446 "TEdgeManager" $null
Fill-RowSolo 447 "TEdgeManager" $null
Fill-RowSolo 448 "TEdgeManager" $null
Fill-RowSolo 449 "TEdgeManager" $null
Fill-RowSolo 450 "TEdgeManager" $null
Fill-RowSolo 451 "TEdgeManager" $null
Fill-RowSolo 452 "TEdgeManager" $null
Fill-RowSolo 453 "TEdgeManager" $null

# --- rows 532 & 535 (Tline) already categorised; just mark Done ---
$ws.Range("B532").Value2 = "Done"
$ws.Range("B535").Value2 = "Done"

# --- row 562 (TPinballTable) flips from "Inp" to "Done" ---
$ws.Range("B562").Value2 = "Done"

# --- brand-new "TPlunger" class: 6 functions done ---
Fill-RowGroup 563 "TPlunger" $null
Fill-RowGroup 564 "TPlunger" $null
Fill-RowGroup 565 "TPlunger" $null
Fill-RowGroup 566 "TPlunger" $null
Fill-RowGroup 567 "TPlunger" $null
Fill-RowGroup 568 "TPlunger" $null

# --- row 606 (TTableLayer) already categorised; just mark Done ---
$ws.Range("B606").Value2 = "Done"

$excel.CutCopyMode = 0

# --- restore the selection/scroll state left by the author ---
$ws.Range("C72").Select()

Write-Host "edits applied"
